# Applies the BAX_cashflow update:
#  - Column widths for C and D narrow to match the rest (15.4)
#  - Row 6 (Change in inventories) and Row 7 (Change in payables and accrued
#    liability) get revised figures for columns B-G
#  - B28 (Capital Stock Change) gets a new numeric value, was blank

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow columns C and D to match the other data columns (~15.4 chars,
# i.e. the same on-screen width Excel uses for the rest of the columns)
$ws.Range("C1").ColumnWidth = 14.714285714285714
$ws.Range("D1").ColumnWidth = 14.714285714285714

# Row 6: Change in inventories
$ws.Range("B6").Value = -208000000.0
$ws.Range("C6").Value = -162000000.0
$ws.Range("D6").Value = -213000000.0
$ws.Range("E6").Value = -166000000.0
$ws.Range("F6").Value = 3000000.0
$ws.Range("G6").Value = 4000000.0

# Row 7: Change in payables and accrued liability
$ws.Range("B7").Value = 151000000.0
$ws.Range("C7").Value = 143000000.0
$ws.Range("D7").Value = 73000000.0
$ws.Range("E7").Value = 39000000.0
$ws.Range("F7").Value = 13000000.0
$ws.Range("G7").Value = -212000000.0

# B28: Capital Stock Change - was blank, now has a value
$ws.Range("B28").Value = -569000000.0
